$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 183.70589
$ws.Range("I39").Value = 120.333336
$ws.Range("J39").Value = 335.8
$ws.Range("K39").Value = 361.000008
$ws.Range("L39").Value = 1007.4
$ws.Range("M39").Value = -65.00000799999998
$ws.Range("N39").Value = -1599.4
$ws.Range("H69").Value = 21069.857
$ws.Range("I69").Value = 32500
$ws.Range("J69").Value = 5829.6665
$ws.Range("K69").Value = 97500
$ws.Range("L69").Value = 17488.9995
$ws.Range("M69").Value = -96626
$ws.Range("N69").Value = -19236.9995
$ws.Range("H72").Value = 21069.857
$ws.Range("I72").Value = 32500
$ws.Range("J72").Value = 5829.6665
$ws.Range("K72").Value = 292500
$ws.Range("L72").Value = 52466.9985
$ws.Range("M72").Value = -288132
$ws.Range("N72").Value = -61202.9985
$ws.Range("H92").Value = 344.625
$ws.Range("I92").Value = 571
$ws.Range("J92").Value = 118.25
$ws.Range("K92").Value = 571
$ws.Range("L92").Value = 118.25
$ws.Range("M92").Value = 677
$ws.Range("N92").Value = -2614.25
$ws.Range("H112").Value = 1129.2
$ws.Range("I112").Value = 100
$ws.Range("J112").Value = 1172.0834
$ws.Range("K112").Value = 300
$ws.Range("L112").Value = 3516.2502
$ws.Range("M112").Value = 808
$ws.Range("N112").Value = -5732.2502
$ws.Range("H137").Value = 501909.4
$ws.Range("I137").Value = 1524.762
$ws.Range("J137").Value = 1815419.1
$ws.Range("K137").Value = 4574.286
$ws.Range("L137").Value = 5446257.300000001
$ws.Range("M137").Value = -2024.286
$ws.Range("N137").Value = -5451357.300000001
$ws.Range("H141").Value = 4342.769
$ws.Range("I141").Value = 3525.7
$ws.Range("J141").Value = 7066.3335
$ws.Range("K141").Value = 10577.1
$ws.Range("L141").Value = 21199.0005
$ws.Range("M141").Value = -5397.099999999999
$ws.Range("N141").Value = -31559.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7644.969
$ws.Range("I32").Value = 4008.6792
$ws.Range("J32").Value = 23705.25
$ws.Range("K32").Value = 4008.6792
$ws.Range("L32").Value = 23705.25
$ws.Range("M32").Value = -3721.6792
$ws.Range("N32").Value = -24279.25
$ws.Range("H45").Value = 12780.091
$ws.Range("I45").Value = 15947.75
$ws.Range("J45").Value = 4333
$ws.Range("K45").Value = 15947.75
$ws.Range("L45").Value = 4333
$ws.Range("M45").Value = -15570.75
$ws.Range("N45").Value = -5087
$ws.Range("H74").Value = 29027.945
$ws.Range("I74").Value = 43124.293
$ws.Range("J74").Value = 3003.923
$ws.Range("K74").Value = 43124.293
$ws.Range("L74").Value = 3003.923
$ws.Range("M74").Value = -42250.293
$ws.Range("N74").Value = -4751.923
$ws.Range("H77").Value = 29027.945
$ws.Range("I77").Value = 43124.293
$ws.Range("J77").Value = 3003.923
$ws.Range("K77").Value = 215621.465
$ws.Range("L77").Value = 15019.615
$ws.Range("M77").Value = -211253.465
$ws.Range("N77").Value = -23755.615
$ws.Range("H102").Value = 75503.87
$ws.Range("I102").Value = 101792.09
$ws.Range("J102").Value = 3211.25
$ws.Range("K102").Value = 101792.09
$ws.Range("L102").Value = 3211.25
$ws.Range("M102").Value = -100170.09
$ws.Range("N102").Value = -6455.25
$ws.Range("H122").Value = 3842.0435
$ws.Range("I122").Value = 3875.5
$ws.Range("J122").Value = 3816.3076
$ws.Range("K122").Value = 11626.5
$ws.Range("L122").Value = 11448.9228
$ws.Range("M122").Value = -9176.5
$ws.Range("N122").Value = -16348.9228
$ws.Range("H132").Value = 1540.2069
$ws.Range("I132").Value = 1519.8
$ws.Range("J132").Value = 1914.3334
$ws.Range("K132").Value = 4559.4
$ws.Range("L132").Value = 5743.0002
$ws.Range("M132").Value = -2029.4
$ws.Range("N132").Value = -10803.0002
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 87816.164
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 87816.164
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 87816.164
$ws.Range("N135").Value = -97956.164

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1233.5555
$ws.Range("I94").Value = 1034.5
$ws.Range("J94").Value = 1631.6666
$ws.Range("K94").Value = 1034.5
$ws.Range("L94").Value = 1631.6666
$ws.Range("M94").Value = -583.5
$ws.Range("N94").Value = -2533.6666
$ws.Range("H134").Value = 1661.7307
$ws.Range("I134").Value = 1532.2273
$ws.Range("J134").Value = 2374
$ws.Range("K134").Value = 4596.6819
$ws.Range("L134").Value = 7122
$ws.Range("M134").Value = -2061.6819
$ws.Range("N134").Value = -12192
$ws.Range("H138").Value = 83757.78
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 83757.78
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 83757.78
$ws.Range("N138").Value = -94037.78

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6804.3335
$ws.Range("I31").Value = 2477.75
$ws.Range("J31").Value = 10265.6
$ws.Range("K31").Value = 2477.75
$ws.Range("L31").Value = 10265.6
$ws.Range("M31").Value = -2182.75
$ws.Range("N31").Value = -10855.6
$ws.Range("H34").Value = 6804.3335
$ws.Range("I34").Value = 2477.75
$ws.Range("J34").Value = 10265.6
$ws.Range("K34").Value = 2477.75
$ws.Range("L34").Value = 10265.6
$ws.Range("M34").Value = -2275.75
$ws.Range("N34").Value = -10669.6
$ws.Range("H121").Value = 38333.668
$ws.Range("I121").Value = 65001
$ws.Range("H132").Value = 1344.4
$ws.Range("I132").Value = 534.8570999999999
$ws.Range("J132").Value = 3233.3333
$ws.Range("K132").Value = 1604.5713
$ws.Range("L132").Value = 9699.999899999999
$ws.Range("M132").Value = 925.4287000000002
$ws.Range("N132").Value = -14759.9999
$ws.Range("H134").Value = 32759.594
$ws.Range("I134").Value = 1378.3572
$ws.Range("J134").Value = 252428.25
$ws.Range("K134").Value = 4135.071599999999
$ws.Range("L134").Value = 757284.75
$ws.Range("M134").Value = -1600.071599999999
$ws.Range("N134").Value = -762354.75
$ws.Range("H140").Value = 73666.336
$ws.Range("I140").Value = 73000
$ws.Range("J140").Value = 73999.5
$ws.Range("K140").Value = 73000
$ws.Range("L140").Value = 73999.5
$ws.Range("M140").Value = -67820
$ws.Range("N140").Value = -84359.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33335586
$ws.Range("I80").Value = 62502150
$ws.Range("J80").Value = 2365.0715
$ws.Range("K80").Value = 62502150
$ws.Range("L80").Value = 2365.0715
$ws.Range("M80").Value = -62501152
$ws.Range("N80").Value = -4361.0715
$ws.Range("H83").Value = 33335586
$ws.Range("I83").Value = 62502150
$ws.Range("J83").Value = 2365.0715
$ws.Range("K83").Value = 312510750
$ws.Range("L83").Value = 11825.3575
$ws.Range("M83").Value = -312505758
$ws.Range("N83").Value = -21809.3575
$ws.Range("H122").Value = 19139.6
$ws.Range("I122").Value = 27067.666
$ws.Range("J122").Value = 7247.5
$ws.Range("K122").Value = 81202.99800000001
$ws.Range("L122").Value = 21742.5
$ws.Range("M122").Value = -78752.99800000001
$ws.Range("N122").Value = -26642.5
$ws.Range("H132").Value = 2859.9575
$ws.Range("I132").Value = 2512.6667
$ws.Range("J132").Value = 3996.5454
$ws.Range("K132").Value = 7538.000100000001
$ws.Range("L132").Value = 11989.6362
$ws.Range("M132").Value = -5008.000100000001
$ws.Range("N132").Value = -17049.6362
$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 42871.75
$ws.Range("I7").Value = 32352.857
$ws.Range("J7").Value = 57598.2
$ws.Range("K7").Value = 32352.857
$ws.Range("L7").Value = 57598.2
$ws.Range("M7").Value = -32240.857
$ws.Range("N7").Value = -57822.2
$ws.Range("H55").Value = 1824.625
$ws.Range("I55").Value = 937.65515
$ws.Range("J55").Value = 4163
$ws.Range("K55").Value = 937.65515
$ws.Range("L55").Value = 4163
$ws.Range("M55").Value = -764.65515
$ws.Range("N55").Value = -4509
$ws.Range("H93").Value = 1932.7222
$ws.Range("I93").Value = 1645.7858
$ws.Range("J93").Value = 2937
$ws.Range("K93").Value = 1645.7858
$ws.Range("L93").Value = 2937
$ws.Range("M93").Value = -397.7858000000001
$ws.Range("N93").Value = -5433
$ws.Range("H126").Value = 42871.75
$ws.Range("I126").Value = 32352.857
$ws.Range("J126").Value = 57598.2
$ws.Range("K126").Value = 97058.571
$ws.Range("L126").Value = 172794.6
$ws.Range("M126").Value = -94588.571
$ws.Range("N126").Value = -177734.6
$ws.Range("H132").Value = 1840.3125
$ws.Range("I132").Value = 1118.9565
$ws.Range("J132").Value = 3683.7778
$ws.Range("K132").Value = 3356.8695
$ws.Range("L132").Value = 11051.3334
$ws.Range("M132").Value = -826.8694999999998
$ws.Range("N132").Value = -16111.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 36428.453
$ws.Range("I94").Value = 25694
$ws.Range("J94").Value = 38813.89
$ws.Range("K94").Value = 25694
$ws.Range("L94").Value = 38813.89
$ws.Range("M94").Value = -24793
$ws.Range("N94").Value = -40615.89
$ws.Range("H98").Value = 8500
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 8500
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 8500
$ws.Range("N98").Value = -14490
$ws.Range("H122").Value = 4064
$ws.Range("I122").Value = 2838.8
$ws.Range("J122").Value = 4829.75
$ws.Range("K122").Value = 8516.400000000001
$ws.Range("L122").Value = 14489.25
$ws.Range("M122").Value = -6066.400000000001
$ws.Range("N122").Value = -19389.25
$ws.Range("H132").Value = 1451029.6
$ws.Range("I132").Value = 1494.4762
$ws.Range("J132").Value = 4833278.5
$ws.Range("K132").Value = 4483.4286
$ws.Range("L132").Value = 14499835.5
$ws.Range("M132").Value = -1953.4286
$ws.Range("N132").Value = -14504895.5
$ws.Range("H136").Value = 4821.55
$ws.Range("I136").Value = 2143
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 6429
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -3879
$ws.Range("N136").Value = -65100
